$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# figures with the latest poll results from the GitHub Actions job.
#
# Column D values are stored as plain text (several look numeric, e.g. "1.000" or
# "0.2580", and some use "." as a thousands separator, e.g. "25.948.68"), so each
# Price cell is forced to Text via NumberFormat before the value is written and the
# number format is reset back to Normal afterwards -- this keeps Excel from silently
# re-interpreting the text as a number (which would drop significant trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.948.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.630.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5232'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2580'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06272'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.49'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07579'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.635.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.421'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.852.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5505'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8008'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.948.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("E19").Value = '  -0.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.674'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '185.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("E22").Value = '  -2.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.108'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1214'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.381'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05897'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.244'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.420'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.388'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.625'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9794'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.382'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.726'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5789'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.69%  '
$ws.Range("E39").Value = '  -1.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8479'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.39%  '
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.036.25'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.675'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.94'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈108'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.004'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.032'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4218'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.86%  '
